# Apply "repull data, push all data, mean calculation" updates to column F (dSF)
# on Sheet1. Only specific rows' F-column values change; everything else is left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = -4
    5  = 1
    6  = 1
    10 = -3
    12 = -1
    17 = -2
    18 = 4
    19 = -4
    23 = -2
    25 = 2
    28 = 0
    30 = -1
    31 = -2
    33 = 1
    43 = 2
    44 = 6
    50 = -1
    56 = -1
    57 = -1
    58 = 4
    60 = -2
    63 = 1
    64 = 6
    65 = -6
    68 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

$wb.Save()
